$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N7").Value = 270
$ws.Range("N9").Value = 270
